$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells for team record columns (AD, AE, AF) on row 1
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the header formatting (bold, centered, bordered) from an existing
# header cell (A1) onto the new header cells
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Restore the values after the paste-special (paste formats only, but
# setting values explicitly again ensures correctness)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the team record (Wins=93, Losses=69, Ties=0) for every data row
for ($r = 2; $r -le 42; $r++) {
    $ws.Cells.Item($r, 30).Value = 93
    $ws.Cells.Item($r, 31).Value = 69
    $ws.Cells.Item($r, 32).Value = 0
}
